# Troca google pela hostinger
#
# The sheet previously listed 5 service-desk tickets (rows 2-6). The edit
# collapses that down to a single placeholder/test row: row 2 is overwritten
# with sample values ("teste" / "10:00" / "11:00" / "01:00" / "Rivaldo" /
# "26/4/2022") and rows 3-6 are deleted outright, shrinking the used range
# from A1:G6 down to A1:G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3-6 entirely (not just their contents) so the sheet's used
# range shrinks back down to A1:G2, matching the target dimension.
$ws.Range("A3:G6").EntireRow.Delete()

# Overwrite the remaining data row (row 2) with the new placeholder values.
$ws.Range("A2").Value = "teste"
$ws.Range("B2").Value = "teste"
$ws.Range("C2").Value = "10:00"
$ws.Range("D2").Value = "11:00"
$ws.Range("E2").Value = "01:00"
$ws.Range("F2").Value = "Rivaldo"
$ws.Range("G2").Value = "26/4/2022"
